# Sprint8.xlsx - "Sprint" sheet: log effort spent fixing bug id25.
#
# The "Implement fix" task (row 6) gets 1 unit of effort logged against
# "Day 5" (column K) -- the bug-fix work that was done that day.
# Every other number on the sheet (U6, K8, U8, K10:T10, the burndown-chart
# feed cells, etc.) is a formula that depends on this cell, so simply
# writing the value and letting Excel recalculate reproduces the rest of
# the sprint table automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint")

# Log 1 unit of effort for the bug fix on Day 5 (column K, row 6 = "Implement fix").
$ws.Range("K6").Value = 1

# Force a recalculation so the dependent totals/remaining-effort formulas
# (U6, K8, U8, K10:T10, ...) pick up the new figure right away.
$excel.Calculate()

# The editor's cursor ended up on K7 (the cell right below the one that was
# just edited) when they were done.
$ws.Range("K7").Select()
